$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-10: Gear (column D) changes from "1-RAP" to "2-RAP"
for ($r = 4; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = "2-RAP"
}

# Row 11: SpecCode (column E) changes from GOBINIG to CALLSAP
$ws.Cells.Item(11, 5).Value = "CALLSAP"

# Row 12: Gear 2-RAP -> 1-RAP, SpecCode ARNOLAT -> TRIGLUC
$ws.Cells.Item(12, 4).Value = "1-RAP"
$ws.Cells.Item(12, 5).Value = "TRIGLUC"

# Row 13: Gear 2-RAP -> 1-RAP, SpecCode PTEOBOV -> GOBINIG
$ws.Cells.Item(13, 4).Value = "1-RAP"
$ws.Cells.Item(13, 5).Value = "GOBINIG"

# New rows 14-17 hold the records that used to be condensed into the old row 14.
# Column C ("Station") must stay text ("4"), matching the rest of the sheet.
$ws.Range("C14:C17").NumberFormat = "@"

# Row 14: new record 1-RAP / ARNOLAT (reuses the row that used to hold 2-RAP/SOLEAEG)
$ws.Cells.Item(14, 1).Value = "SOLEMON2024"
$ws.Cells.Item(14, 2).Value = "ITA17"
$ws.Cells.Item(14, 3).Value = "4"
$ws.Cells.Item(14, 4).Value = "1-RAP"
$ws.Cells.Item(14, 5).Value = "ARNOLAT"
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = -1
$ws.Cells.Item(14, 8).Value = "SIMRANDO"
$ws.Cells.Item(14, 9).Value = "Y"

# New row 15: 1-RAP / SEPIOFF
$ws.Cells.Item(15, 1).Value = "SOLEMON2024"
$ws.Cells.Item(15, 2).Value = "ITA17"
$ws.Cells.Item(15, 3).Value = "4"
$ws.Cells.Item(15, 4).Value = "1-RAP"
$ws.Cells.Item(15, 5).Value = "SEPIOFF"
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = -1
$ws.Cells.Item(15, 8).Value = "SIMRANDO"
$ws.Cells.Item(15, 9).Value = "Y"

# New row 16: 2-RAP / PTEOBOV
$ws.Cells.Item(16, 1).Value = "SOLEMON2024"
$ws.Cells.Item(16, 2).Value = "ITA17"
$ws.Cells.Item(16, 3).Value = "4"
$ws.Cells.Item(16, 4).Value = "2-RAP"
$ws.Cells.Item(16, 5).Value = "PTEOBOV"
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = -1
$ws.Cells.Item(16, 8).Value = "SIMRANDO"
$ws.Cells.Item(16, 9).Value = "Y"

# New row 17: 2-RAP / SOLEAEG (the data that used to live in the old row 14)
$ws.Cells.Item(17, 1).Value = "SOLEMON2024"
$ws.Cells.Item(17, 2).Value = "ITA17"
$ws.Cells.Item(17, 3).Value = "4"
$ws.Cells.Item(17, 4).Value = "2-RAP"
$ws.Cells.Item(17, 5).Value = "SOLEAEG"
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = -1
$ws.Cells.Item(17, 8).Value = "SIMRANDO"
$ws.Cells.Item(17, 9).Value = "Y"
